$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 83351600
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 250050000
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 250050000
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -250051248
$ws.Range("H65").Value = 83351600
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 250050000
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 1250250000
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -1250256240
$ws.Range("H113").Value = 71431384
$ws.Range("I113").Value = 250001000
$ws.Range("J113").Value = 3536.8
$ws.Range("K113").Value = 250001000
$ws.Range("L113").Value = 3536.8
$ws.Range("M113").Value = -249997746
$ws.Range("N113").Value = -10044.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2628.5454
$ws.Range("I45").Value = 1797.174
$ws.Range("J45").Value = 4540.7
$ws.Range("K45").Value = 1797.174
$ws.Range("L45").Value = 4540.7
$ws.Range("M45").Value = -1420.174
$ws.Range("N45").Value = -5294.7
$ws.Range("H61").Value = 1272
$ws.Range("I61").Value = 1395.6666
$ws.Range("J61").Value = 530
$ws.Range("K61").Value = 1395.6666
$ws.Range("L61").Value = 530
$ws.Range("M61").Value = -1183.6666
$ws.Range("N61").Value = -954
$ws.Range("H110").Value = 7751.8887
$ws.Range("I110").Value = 7740.364
$ws.Range("J110").Value = 7802.6
$ws.Range("K110").Value = 7740.364
$ws.Range("L110").Value = 7802.6
$ws.Range("M110").Value = -5695.364
$ws.Range("N110").Value = -11892.6
$ws.Range("H136").Value = 1272
$ws.Range("I136").Value = 1395.6666
$ws.Range("J136").Value = 530
$ws.Range("K136").Value = 4186.9998
$ws.Range("L136").Value = 1590
$ws.Range("M136").Value = -1636.9998
$ws.Range("N136").Value = -6690

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 159020.31
$ws.Range("I134").Value = 231825.16
$ws.Range("J134").Value = 1276.5
$ws.Range("K134").Value = 695475.48
$ws.Range("L134").Value = 3829.5
$ws.Range("M134").Value = -692940.48
$ws.Range("N134").Value = -8899.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1739.4
$ws.Range("I99").Value = 1650
$ws.Range("J99").Value = 1873.5
$ws.Range("K99").Value = 1650
$ws.Range("L99").Value = 1873.5
$ws.Range("M99").Value = -152
$ws.Range("N99").Value = -4869.5
$ws.Range("H126").Value = 1739.4
$ws.Range("I126").Value = 1650
$ws.Range("J126").Value = 1873.5
$ws.Range("K126").Value = 4950
$ws.Range("L126").Value = 5620.5
$ws.Range("M126").Value = -2480
$ws.Range("N126").Value = -10560.5
$ws.Range("H134").Value = 3571.7878
$ws.Range("I134").Value = 3771.3965
$ws.Range("K134").Value = 11314.1895
$ws.Range("M134").Value = -8779.1895
$ws.Range("H135").Value = 54995.645
$ws.Range("J135").Value = 54995.645
$ws.Range("L135").Value = 54995.645
$ws.Range("N135").Value = -65135.645

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 4213.278
$ws.Range("J123").Value = 4377.8823
$ws.Range("L123").Value = 13133.6469
$ws.Range("N123").Value = -18033.6469
$ws.Range("H131").Value = 914.6774
$ws.Range("J131").Value = 976.4815
$ws.Range("L131").Value = 2929.4445
$ws.Range("N131").Value = -13009.4445

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1867.8182
$ws.Range("I7").Value = 1482.8889
$ws.Range("J7").Value = 3600
$ws.Range("K7").Value = 1482.8889
$ws.Range("L7").Value = 3600
$ws.Range("M7").Value = -1370.8889
$ws.Range("N7").Value = -3824
$ws.Range("H61").Value = 1773.2307
$ws.Range("I61").Value = 1504.8182
$ws.Range("K61").Value = 1504.8182
$ws.Range("M61").Value = -1302.8182
$ws.Range("H68").Value = 1894.7368
$ws.Range("I68").Value = 1782.3529
$ws.Range("J68").Value = 2850
$ws.Range("K68").Value = 1782.3529
$ws.Range("L68").Value = 2850
$ws.Range("M68").Value = -1033.3529
$ws.Range("N68").Value = -4348
$ws.Range("H71").Value = 1894.7368
$ws.Range("I71").Value = 1782.3529
$ws.Range("J71").Value = 2850
$ws.Range("K71").Value = 8911.764500000001
$ws.Range("L71").Value = 14250
$ws.Range("M71").Value = -5167.764500000001
$ws.Range("N71").Value = -21738
$ws.Range("H113").Value = 1773.2307
$ws.Range("I113").Value = 1504.8182
$ws.Range("K113").Value = 1504.8182
$ws.Range("M113").Value = 665.1818000000001
$ws.Range("H118").Value = 39888
$ws.Range("J118").Value = 39888
$ws.Range("L118").Value = 39888
$ws.Range("N118").Value = -43202
$ws.Range("H126").Value = 1867.8182
$ws.Range("I126").Value = 1482.8889
$ws.Range("J126").Value = 3600
$ws.Range("K126").Value = 4448.6667
$ws.Range("L126").Value = 10800
$ws.Range("M126").Value = -1978.6667
$ws.Range("N126").Value = -15740
$ws.Range("H132").Value = 2511.9375
$ws.Range("I132").Value = 2269.95
$ws.Range("J132").Value = 2915.25
$ws.Range("K132").Value = 6809.849999999999
$ws.Range("L132").Value = 8745.75
$ws.Range("M132").Value = -4279.849999999999
$ws.Range("N132").Value = -13805.75
$ws.Range("H136").Value = 1915.6957
$ws.Range("I136").Value = 1671.4865
$ws.Range("J136").Value = 2919.6667
$ws.Range("K136").Value = 5014.4595
$ws.Range("L136").Value = 8759.000100000001
$ws.Range("M136").Value = -2464.4595
$ws.Range("N136").Value = -13859.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1798.6945
$ws.Range("I126").Value = 1409.0358
$ws.Range("J126").Value = 3162.5
$ws.Range("K126").Value = 4227.107400000001
$ws.Range("L126").Value = 9487.5
$ws.Range("M126").Value = -1757.107400000001
$ws.Range("N126").Value = -14427.5
$ws.Range("H136").Value = 1695.4565
$ws.Range("I136").Value = 1793.2703
$ws.Range("J136").Value = 1293.3334
$ws.Range("K136").Value = 5379.810899999999
$ws.Range("L136").Value = 3880.0002
$ws.Range("M136").Value = -2829.810899999999
$ws.Range("N136").Value = -8980.0002

Write-Host "Applied all profit-sheet updates"
